$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Value = 3.89
$ws.Range("D16").Value = 0.77
$ws.Range("D36").Value = 3.76
$ws.Range("D37").Value = 3.73
$ws.Range("D58").Value = 3.39
